$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 1184
$ws1.Range("F11").Value = 3075
$ws1.Range("F15").Value = 624
$ws1.Range("F19").Value = 1156
$ws1.Range("F20").Value = 1156
$ws1.Range("F25").Value = 262
$ws1.Range("F30").Value = 102
$ws1.Range("F35").Value = 543
$ws1.Range("F36").Value = 282

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F38").Value = 453

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 1184
$ws4.Range("F16").Value = 3075
$ws4.Range("F25").Value = 1156
$ws4.Range("F26").Value = 1156
$ws4.Range("F32").Value = 262
$ws4.Range("F39").Value = 102
$ws4.Range("F46").Value = 543
$ws4.Range("F47").Value = 453
$ws4.Range("F48").Value = 282
